$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.858.53"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.630.57"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.792.80"
$ws.Range("E12").Value = "  +9.77%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.855.18"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "0.0₃0754"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "25.871.06"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").Value = "1.136.68"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "1.764.87"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0113"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0525"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
